# AP API Test Cases
# Reproduces the authoring session captured in the commit:
#   - PAYTO gets three new test rows (Cash / Bank Card paytypes) and becomes
#     the active/selected sheet (selection on row 3).
#   - ObjectName is no longer the selected/active sheet.
#   - A handful of cells that carried a redundant "applyFill" style (with no
#     actual fill) are reset back to the workbook's default (unstyled) look.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Strip the stray no-op fill style off the cells that still reference it.
#    (numFmtId=0 / fontId=0 / fillId=0 "applyFill" xf -- visually identical
#    to the default style, so clearing the interior pattern drops the
#    reference entirely.)
# ---------------------------------------------------------------------------

$noStyleTargets = @{
    1  = @("A2", "H2")
    4  = @("A3")
    7  = @(
        "A2","B2","E2",
        "A3","B3","E3",
        "A4","B4","E4",
        "A5","B5","E5",
        "A6","B6","E6",
        "A7","B7","E7",
        "A8","B8","E8",
        "A9","B9","E9",
        "A10","B10","E10",
        "A11","B11","E11",
        "A12","B12","E12",
        "A13","B13","E13",
        "A14","B14","E14"
    )
    8  = @(
        "A2","B2","E2",
        "A3","B3","E3",
        "A4","B4","E4",
        "A5","B5","E5",
        "A6","B6","E6",
        "A7","B7","E7",
        "A8","B8","E8",
        "A9","B9","E9",
        "A10","B10","E10",
        "A11","B11","E11",
        "A12","B12","E12",
        "A13","B13","E13",
        "A14","B14","E14"
    )
    9  = @("A2", "B2", "E2")
    10 = @("A2", "B2", "C2", "E2", "F2", "G2")
}

foreach ($sheetIndex in $noStyleTargets.Keys) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    foreach ($addr in $noStyleTargets[$sheetIndex]) {
        $ws.Range($addr).Interior.Pattern = -4142   # xlNone
    }
}

# ---------------------------------------------------------------------------
# 2. PAYTO: add the three new AP API test rows.
# ---------------------------------------------------------------------------

$payto = $wb.Worksheets.Item(9)

$payto.Range("A3").Value = "Mountain Manufacturing (100)"
$payto.Range("B3").Value = "6655 (BC)"
$payto.Range("C3").Value = "EFT"
$payto.Range("D3").Value = 2
$payto.Range("E3").Value = "SB-24808 (1103)"
$payto.Range("F3").Value = 230
$payto.Range("G3").Value = "Open"
$payto.Range("H3").Value = "Invoice"
$payto.Range("I3").Value = $true
$payto.Range("J3").Value = $false

# Row 5 is written before row 4 so the new shared strings land in the same
# order as the workbook being reproduced ("Cash" before "Bank Card").
$payto.Range("A5").Value = "Mountain Manufacturing (100)"
$payto.Range("B5").Value = "6655 (BC)"
$payto.Range("C5").Value = "Cash"
$payto.Range("D5").Value = 3
$payto.Range("E5").Value = "SB-24808 (1103)"
$payto.Range("F5").Value = 340
$payto.Range("G5").Value = "Open"
$payto.Range("H5").Value = "Invoice"
$payto.Range("I5").Value = $true
$payto.Range("J5").Value = $false

$payto.Range("A4").Value = "Mountain Manufacturing (100)"
$payto.Range("B4").Value = "6655 (BC)"
$payto.Range("C4").Value = "Bank Card"
$payto.Range("D4").Value = 3
$payto.Range("E4").Value = "SB-24808 (1103)"
$payto.Range("F4").Value = 340
$payto.Range("G4").Value = "Open"
$payto.Range("H4").Value = "Invoice"
$payto.Range("I4").Value = $true
$payto.Range("J4").Value = $false

# ---------------------------------------------------------------------------
# 3. Make PAYTO the active sheet with row 3 selected; ObjectName (previously
#    active) loses the selection.
# ---------------------------------------------------------------------------

$payto.Select()
$payto.Range("A3:XFD3").Select()
